# Update "想去人数" (people interested) figures for two events that
# appear on both the "展览" sheet and the "全部类型" sheet.
#   - 南宁·2024良牙动漫冬季盛典（冬典） : 9100 -> 9112 (cell F2)
#   - 南宁·草莓动漫节                   : 469  -> 470  (cell F4)

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 9112
    $ws.Range("F4").Value = 470
}
